$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts C:J -> D:K)
$ws.Columns("C").Insert()

# New column header
$ws.Range("C1").Value = "M_PL"

# New column data values
$ws.Range("C2").Value = 141492655254
$ws.Range("C3").Value = 1017808846331
$ws.Range("C4").Value = 41087099249
$ws.Range("C5").Value = 51692957042
$ws.Range("C6").Value = 959424197928
$ws.Range("C7").Value = 42657873476
$ws.Range("C8").Value = 91572151625
$ws.Range("C9").Value = 1076750542353
$ws.Range("C10").Value = 208898496050
$ws.Range("C11").Value = 867852046303
